# Replace the leftover "<type>text</type>" template placeholders that sit
# inside the summary-table cells with real Word runs (w:r/w:t), matching the
# move to directly emitting TextRun content instead of the old placeholder
# marker.

$d = $word.ActiveDocument

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Set-CellRun($Table, $Row, $Col, $Text, $RunProps) {

    $cell = $Table.Cell($Row, $Col)
    $para = $cell.Range.Paragraphs.Item(1)

    $xml = '<w:p ' + $wNs + '>' +
             '<w:pPr><w:spacing w:before="40" w:after="40"/></w:pPr>' +
             '<w:r>' +
               $RunProps +
               '<w:t xml:space="preserve">' + $Text + '</w:t>' +
             '</w:r>' +
           '</w:p>'

    [void]$para.Range.InsertXML($xml)
}

# Shared run-formatting fragments used throughout the two summary tables.
$labelProps = '<w:rPr>' +
                '<w:rFonts w:ascii="Arial" w:cs="Arial" w:eastAsia="Arial" w:hAnsi="Arial"/>' +
                '<w:color w:val="6b7280"/>' +
                '<w:sz w:val="20"/>' +
                '<w:szCs w:val="20"/>' +
              '</w:rPr>'

$valueProps = '<w:rPr>' +
                '<w:rFonts w:ascii="Arial" w:cs="Arial" w:eastAsia="Arial" w:hAnsi="Arial"/>' +
                '<w:b/>' +
                '<w:bCs/>' +
                '<w:sz w:val="22"/>' +
                '<w:szCs w:val="22"/>' +
              '</w:rPr>'

# First summary table: 4 columns x 2 rows (Name/Date of evaluation, Date of
# birth/Evaluation ID).
$t1 = $d.Tables.Item(2)

Set-CellRun $t1 1 1 "NAME OF APPLICANT" $labelProps
Set-CellRun $t1 1 2 "[first name] [last name]" $valueProps
Set-CellRun $t1 1 3 "DATE OF EVALUATION" $labelProps
Set-CellRun $t1 1 4 "January 16, 2026" $valueProps

Set-CellRun $t1 2 1 "DATE OF BIRTH" $labelProps
Set-CellRun $t1 2 2 "N/A" $valueProps
Set-CellRun $t1 2 3 "EVALUATION ID" $labelProps
Set-CellRun $t1 2 4 "LA-20260116-111" $valueProps

# Second summary table: 2 columns x 2 rows (Purpose of evaluation, Country of
# education).
$t2 = $d.Tables.Item(3)

Set-CellRun $t2 1 1 "PURPOSE OF EVALUATION" $labelProps
Set-CellRun $t2 1 2 "N/A" $valueProps

Set-CellRun $t2 2 1 "COUNTRY OF EDUCATION" $labelProps
Set-CellRun $t2 2 2 "Poland" $valueProps

Write-Output "done"
